$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1018
$ws.Range("F7").Value = 2576
$ws.Range("F9").Value = 1258
$ws.Range("F10").Value = 911
$ws.Range("F11").Value = 609
$ws.Range("F13").Value = 1148
$ws.Range("F17").Value = 734
$ws.Range("F18").Value = 783
$ws.Range("F19").Value = 206
$ws.Range("F21").Value = 1123
$ws.Range("F22").Value = 96
$ws.Range("F23").Value = 615
$ws.Range("F24").Value = 597
$ws.Range("F25").Value = 218
$ws.Range("F26").Value = 304
$ws.Range("F27").Value = 305
$ws.Range("F28").Value = 684
$ws.Range("F29").Value = 488
$ws.Range("F30").Value = 4503
$ws.Range("F35").Value = 157
$ws.Range("F36").Value = 1615
$ws.Range("F38").Value = 55
$ws.Range("F39").Value = 441
$ws.Range("F43").Value = 72
$ws.Range("F45").Value = 133
$ws.Range("F47").Value = 111
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 99
$ws.Range("F9").Value = 23
$ws.Range("F14").Value = 23
$ws.Range("F16").Value = 31
$ws.Range("F17").Value = 194
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2272
$ws.Range("F3").Value = 730
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2272
$ws.Range("F3").Value = 730
$ws.Range("F7").Value = 1018
$ws.Range("F8").Value = 2576
$ws.Range("F10").Value = 1258
$ws.Range("F11").Value = 911
$ws.Range("F12").Value = 609
$ws.Range("F14").Value = 1148
$ws.Range("F18").Value = 734
$ws.Range("F19").Value = 783
$ws.Range("F20").Value = 206
$ws.Range("F22").Value = 1123
$ws.Range("F23").Value = 99
$ws.Range("F24").Value = 96
$ws.Range("F25").Value = 615
$ws.Range("F26").Value = 597
$ws.Range("F27").Value = 218
$ws.Range("F28").Value = 305
$ws.Range("F30").Value = 488
$ws.Range("F31").Value = 4503
$ws.Range("F36").Value = 157
$ws.Range("F37").Value = 1615
$ws.Range("F38").Value = 441
$ws.Range("F44").Value = 72
$ws.Range("F46").Value = 133
$ws.Range("F48").Value = 111
